$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Version: 5.0.0 -> 6.0.0
$ws.Range("B3").Value = "6.0.0"

# Date: refreshed publish date
$ws.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher now has a value
$ws.Range("B9").Value = "Alvearie Team"

# The two "Contact" rows (10 & 11) are replaced by a single "Jurisdiction" row
$ws.Range("A10").Value = "Jurisdiction"
$ws.Range("B10").Value = "United States of America"

# Remove the now-redundant second Contact row; everything below shifts up one row
$ws.Rows.Item(11).Delete()

# "Case Sensitive" (now row 14) gets a literal text value "true"
# (use a formula + paste-special-values so Excel stores it as text, not a boolean)
$ws.Range("B14").Formula = "=""true"""
$ws.Range("B14").Copy()
$ws.Range("B14").PasteSpecial(-4163)
